$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 - Gas any
$ws.Range("C2").Value = 1523
$ws.Range("D2").Value = 460
$ws.Range("E2").Value = 9786
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0.1346715005747635
$ws.Range("H2").Value = 0.1285037521066134
$ws.Range("I2").Value = 0.1410873554787151
$ws.Range("J2").Value = 0.7680282400403429
$ws.Range("K2").Value = 0.7489430154271794
$ws.Range("L2").Value = 0.7860770261975569

# Row 3 - ABG threshold
$ws.Range("C3").Value = 971
$ws.Range("D3").Value = 1012
$ws.Range("E3").Value = 6299
$ws.Range("F3").Value = 3487
$ws.Range("G3").Value = 0.1335625859697387
$ws.Range("H3").Value = 0.1259360465649205
$ws.Range("I3").Value = 0.1415761709939392
$ws.Range("J3").Value = 0.4896621280887544
$ws.Range("K3").Value = 0.4677012993473234
$ws.Range("L3").Value = 0.5116629323491825

# Row 4 - VBG threshold
$ws.Range("C4").Value = 1300
$ws.Range("D4").Value = 683
$ws.Range("E4").Value = 4944
$ws.Range("F4").Value = 4842
$ws.Range("G4").Value = 0.2081998718770019
$ws.Range("H4").Value = 0.1983099660476946
$ws.Range("I4").Value = 0.2184486018485176
$ws.Range("J4").Value = 0.6555723651033787
$ws.Range("K4").Value = 0.6343752062961309
$ws.Range("L4").Value = 0.6761679411003531

# Row 5 - PCO2 OTHER threshold
$ws.Range("C5").Value = 330
$ws.Range("D5").Value = 1653
$ws.Range("E5").Value = 1020
$ws.Range("F5").Value = 8766
$ws.Range("G5").Value = 0.2444444444444444
$ws.Range("H5").Value = 0.2222658923016254
$ws.Range("I5").Value = 0.2680732493279668
$ws.Range("J5").Value = 0.1664145234493192
$ws.Range("K5").Value = 0.1506696927182576
$ws.Range("L5").Value = 0.1834492959275313
